$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly populated "Q" column (year 2018 data) values that were
# previously blank. Existing cell styles are preserved automatically since
# only the Value is being set.
$ws.Range("Q6").Value  = 0.04416600769365872
$ws.Range("Q8").Value  = 601820.3
$ws.Range("Q10").Value = 0.36185407133694547
$ws.Range("Q12").Value = 20892.4
$ws.Range("Q14").Value = 0.274087106792226
$ws.Range("Q16").Value = 63884.8
$ws.Range("Q18").Value = 0.00046658384803364067
$ws.Range("Q20").Value = 85729.5
$ws.Range("Q22").Value = 0.08603265705379398
$ws.Range("Q24").Value = 16970.3
$ws.Range("Q26").Value = "-"
$ws.Range("Q28").Value = 47183.5
$ws.Range("Q30").Value = "-"
$ws.Range("Q32").Value = 17405.3
$ws.Range("Q34").Value = 0.000011900270969169968
$ws.Range("Q36").Value = 84031.7
$ws.Range("Q38").Value = "-"
$ws.Range("Q40").Value = 231841.7
$ws.Range("Q42").Value = "-"
$ws.Range("Q44").Value = 33881.1

# Update the view: scroll back so column A is the left-most visible column
# again, and move the active selection to O52.
$ws.Activate()
$ws.Range("O52").Select()
